$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New columns J (=H*100), K (=I*100), L (relative error) for rows 3..21 ---
# Row 3 entered individually (becomes the "master" of a shared-formula block for rows 4-21)
$ws.Range("J3").Formula = "=H3*100"
$ws.Range("K3").Formula = "=I3*100"
$ws.Range("L3").Formula = "=(K3-J3)/J3"

# Fill J4:K21 in one shot (mirrors dragging the two-column block down) and L4:L21 separately,
# so the engine groups them the same way Excel would (two shared-formula groups).
$ws.Range("J4:K21").Formula = "=H4*100"
$ws.Range("L4:L21").Formula = "=(K4-J4)/J4"

# --- Number formatting ---
$ws.Range("K3:K21").NumberFormat = "0.00"
$ws.Range("L3:L21").Style = "Percent"

# --- Summary cell L24: mean absolute percentage error (legacy CSE array formula) ---
$ws.Range("L24").FormulaArray = "=SUM(ABS(L3:L21))/COUNT(L3:L21)"
$ws.Range("L24").NumberFormat = "0.00%"

# --- Column width for L ---
$ws.Columns.Item("L").ColumnWidth = 20.28515625

# --- Selection / view ---
$ws.Range("J3:L21").Select()

# --- Chart title & series names on the second chart (distance vs RSSI comparison) ---
$chartObjs = $ws.ChartObjects()
$co = $chartObjs.Item(2)
$chart = $co.Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "ESP32 - Obstacles"
$series = $chart.SeriesCollection()
$series.Item(1).Name = "Real RSSI"
$series.Item(2).Name = "Computed RSSI"
